$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a handful of values with refreshed (re-concatenated) figures.
$ws.Range("P57").Value = 549757.056
$ws.Range("P59").Value = 537946.048
$ws.Range("X59").Value = 917802.944
$ws.Range("AN59").Value = 215652.048
$ws.Range("P60").Value = -369773.024
$ws.Range("T60").Value = -744506.048
$ws.Range("X60").Value = -674940.032
$ws.Range("AB60").Value = -839934.016
$ws.Range("AF60").Value = -875300.992
$ws.Range("AN60").Value = -315561.024
$ws.Range("H61").Value = 85892.008
$ws.Range("L61").Value = 49034.992
$ws.Range("X61").Value = 242862.96
$ws.Range("AN61").Value = -99908.984
$ws.Range("T63").Value = -15757
$ws.Range("AF68").Value = 10855
$ws.Range("AN68").Value = -7327
$ws.Range("AB70").Value = -60853
$ws.Range("AF70").Value = -71173.992
$ws.Range("AN70").Value = -30299.992
$ws.Range("L74").Value = -4258.008
$ws.Range("T74").Value = 131582.024
$ws.Range("X74").Value = 98739.992
$ws.Range("AF74").Value = -432664
$ws.Range("T80").Value = 97637.992
$ws.Range("X80").Value = 75122.008
$ws.Range("AB80").Value = -29589
$ws.Range("AJ80").Value = -533837.024
$ws.Range("AN80").Value = -440086.048

# Rows that had placeholder zeros now fall outside the source balance
# sheets reported period -- clear them to blank instead of 0.
$ws.Range("T57:AP57").ClearContents()
$ws.Range("T58:AP58").ClearContents()
$ws.Range("B64:C64").ClearContents()
$ws.Range("E64:T64").ClearContents()
$ws.Range("T71:AP71").ClearContents()
$ws.Range("T72:AP72").ClearContents()
$ws.Range("T73:AP73").ClearContents()
$ws.Range("T77:AP77").ClearContents()
$ws.Range("T78:AP78").ClearContents()
